# ND01.xlsx: "Merge back T2A sheet in the test files"
#
# Adds a new "T2A" worksheet (a near-duplicate of the "WMT_Extract" sheet's
# layout, but with mostly-zeroed numeric data) as the last tab, and makes it
# the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet as the very last tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "T2A"

# Reference cells (already present elsewhere in the workbook) whose cell
# style we reuse, so no redundant/duplicate style or font entries are
# created - only the one genuinely new numFmt/font combination needed for
# the date column gets added.
$srcWmt = $wb.Worksheets.Item("WMT_Extract")
$srcGS  = $wb.Worksheets.Item("GS")
$srcArms = $wb.Worksheets.Item("ARMS")

# style s="2"  (Arial 13, black) - used on the header row
$srcWmt.Range("Y1").Copy() | Out-Null
$ws.Range("A1:AO1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# style s="19" (Calibri 12, black) - used on the data rows
$srcGS.Range("E2").Copy() | Out-Null
$ws.Range("A2:AN4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# style s="16" (Calibri 12, black, numFmtId 14) as a starting point for the
# AO (Datestamp) column, then switch its number format to the m/d/yy h:mm
# one actually used by the source data - this creates the single new
# numFmtId=22/fontId=8 style combination.
$srcArms.Range("A2").Copy() | Out-Null
$ws.Range("AO2:AO4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("AO2:AO4").NumberFormat = "m/d/yy h:mm"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Header row (A1:AO1).
# ---------------------------------------------------------------------
$header = @(
    "Trust", "Region_Desc", "Region_Code", "Ldu_Desc", "Ldu_Code",
    "Team_Desc", "Team_Code", "OM_Surname", "OM_Forename", "OM_Grade_Code",
    "OM_Key", "CommTier0", "CommTierD2", "CommTier1a", "CommTier",
    "CommTierD1", "CommTierC2", "CommTier3a", "CommTierC1", "CommTierB2",
    "CommTierB1", "CommTierA", "LicenceTier0", "LicenceTierD2",
    "LicenceTierD1", "LicenceTierC2", "LicenceTierC1", "LicenceTierB2",
    "LicenceTierB1", "LicenceTierA", "CustTier0", "CustTierD2",
    "CustTierD1", "CustTierC2", "CustTierC1", "CustTierB2", "CustTierB1",
    "CustTierA", "ComIn1st16Weeks", "LicIn1st16Weeks", "Datestamp"
)
for ($i = 0; $i -lt $header.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $header[$i]
}
$ws.Rows.Item(1).RowHeight = 17

# ---------------------------------------------------------------------
# 3. Data rows (A2:J4) - text columns.
# ---------------------------------------------------------------------
$row2 = @("Farringdon", "London", "ND01", "KainosLDU", "KNS", "WMT Team", "WMT", "Swann", "Tom", "C")
$row3 = @("Farringdon", "London", "ND01", "KainosLDU", "KNS", "WMT Team", "WMT", "Wright", "Andy", "Z")
$row4 = @("Farringdon", "London", "ND01", "Jonahs LDU", "JLDU", "WMT Team", "JWMT", "Smith", "Jonah", "C")

$textRows = @{ 2 = $row2; 3 = $row3; 4 = $row4 }
foreach ($r in $textRows.Keys) {
    $vals = $textRows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value2 = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# 4. Data rows - numeric columns K:AN (case ref, count, then all zeros).
# ---------------------------------------------------------------------
$caseRef = @{ 2 = 1001; 3 = 1002; 4 = 1003 }
$countVal = @{ 2 = 1; 3 = 10; 4 = 5 }

foreach ($r in 2, 3, 4) {
    $ws.Cells.Item($r, 11).Value2 = $caseRef[$r]   # column K
    $ws.Cells.Item($r, 12).Value2 = $countVal[$r]  # column L
    for ($c = 13; $c -le 40; $c++) {                # columns M:AN
        $ws.Cells.Item($r, $c).Value2 = 0
    }
}

# ---------------------------------------------------------------------
# 5. Data rows - AO (Datestamp) column.
# ---------------------------------------------------------------------
foreach ($r in 2, 3, 4) {
    $ws.Cells.Item($r, 41).Value2 = 42795.628472222219
}

# ---------------------------------------------------------------------
# 6. Sheet view: select the full used range, no special active cell.
# ---------------------------------------------------------------------
$ws.Range("A1:AO4").Select()

$wb.Save()
